$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename labels: "Te" -> "Tevap" and "pe(5K overheated) / Pa" -> "pevap / Pa" ---
# (Header row for each of the three data blocks on the sheet)
$ws.Range("A3").Value = "Tevap"
$ws.Range("E3").Value = "pevap / Pa"
$ws.Range("M3").Value = "pevap / Pa"

$ws.Range("E15").Value = "pevap / Pa"
$ws.Range("M15").Value = "pevap / Pa"

$ws.Range("E27").Value = "pevap / Pa"

# --- Update input values (dependent formulas recalc automatically) ---
$ws.Range("A4").Value = 0
$ws.Range("I16").Value = 15
$ws.Range("A28").Value = 20

# --- Move / resize the embedded chart ---
$co = $ws.ChartObjects(1)
$co.Left = 569.8544134473425
$co.Top = 88.12496062992126
$co.Width = 683.7295709276575
$co.Height = 338.62503937007875

# --- Update the active selection ---
[void]$ws.Range("E27").Select()
